# Update the "想去人数" (want-to-go count) figures in the 展览 and 全部类型
# sheets to reflect the latest scrape (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" and sheet "全部类型" share the same F-column updates for the
# first six events; the remaining four events live on different rows because
# sheet "全部类型" has two extra rows above them.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F5").Value  = 15588
    $ws.Range("F9").Value  = 15405
    $ws.Range("F11").Value = 9007
    $ws.Range("F12").Value = 379
    $ws.Range("F25").Value = 1110
    $ws.Range("F29").Value = 84
}

# "展览" sheet specific rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F35").Value = 251
$ws1.Range("F36").Value = 319
$ws1.Range("F38").Value = 117
$ws1.Range("F39").Value = 5529

# "全部类型" sheet specific rows (offset by 2 vs. "展览")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F37").Value = 251
$ws4.Range("F38").Value = 319
$ws4.Range("F40").Value = 117
$ws4.Range("F41").Value = 5529
